$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Apply theme-colored fill style to A22 / B22 ---
$ws.Range("A22:B22").Interior.ThemeColor = 8

# --- New "OBSERVAÇÕES" note block (E35:F38) ---
$ws.Range("E35").Value2 = "OBSERVAÇÕES"
$ws.Range("E36").Value2 = "Se for montar com 32K de RAM alta, reduzir um"
$ws.Range("E38").Value2 = "diodos 1N4148 (comprar somente 22)"
$ws.Range("E37").Value2 = "resistor de 1K (comprar somente 4) e reduzir dois"

$ws.Range("E35:F38").Interior.ThemeColor = 8

# --- A42 / B42 fill ---
$ws.Range("A42:B42").Interior.ThemeColor = 8

# --- sheet view ---
$ws.Application.ActiveWindow.ScrollRow = 19
$ws.Range("E38").Select()
